# Daily Slovakia COVID stats update - "Updated: st 27. 05. 2021"
# Applies revised AgTests (F) / AgPosit (G) figures for existing rows and
# appends the new day's record (row 448, date 44342 = 2021-05-26).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to previously reported AgTests (col F) / AgPosit (col G) values ---
$ws.Cells.Item(344, 6).Value2 = 135987
$ws.Cells.Item(344, 7).Value2 = 2485
$ws.Cells.Item(390, 6).Value2 = 219890
$ws.Cells.Item(395, 6).Value2 = 752175
$ws.Cells.Item(395, 7).Value2 = 1948
$ws.Cells.Item(396, 6).Value2 = 166492
$ws.Cells.Item(396, 7).Value2 = 551
$ws.Cells.Item(398, 6).Value2 = 298542
$ws.Cells.Item(399, 6).Value2 = 200445
$ws.Cells.Item(399, 7).Value2 = 968
$ws.Cells.Item(400, 6).Value2 = 149057
$ws.Cells.Item(400, 7).Value2 = 764
$ws.Cells.Item(403, 6).Value2 = 353890
$ws.Cells.Item(404, 6).Value2 = 223951
$ws.Cells.Item(407, 6).Value2 = 158090
$ws.Cells.Item(407, 7).Value2 = 673
$ws.Cells.Item(409, 6).Value2 = 708158
$ws.Cells.Item(410, 6).Value2 = 364460
$ws.Cells.Item(417, 6).Value2 = 342466
$ws.Cells.Item(417, 7).Value2 = 589
$ws.Cells.Item(418, 6).Value2 = 202135
$ws.Cells.Item(418, 7).Value2 = 700
$ws.Cells.Item(419, 6).Value2 = 149297
$ws.Cells.Item(420, 6).Value2 = 138714
$ws.Cells.Item(423, 6).Value2 = 439989
$ws.Cells.Item(424, 6).Value2 = 265658
$ws.Cells.Item(428, 6).Value2 = 102279
$ws.Cells.Item(429, 6).Value2 = 177982
$ws.Cells.Item(430, 6).Value2 = 175314
$ws.Cells.Item(431, 6).Value2 = 170854
$ws.Cells.Item(433, 6).Value2 = 86040
$ws.Cells.Item(434, 6).Value2 = 79560
$ws.Cells.Item(436, 6).Value2 = 145158
$ws.Cells.Item(443, 6).Value2 = 105604
$ws.Cells.Item(444, 6).Value2 = 102797
$ws.Cells.Item(445, 6).Value2 = 83903
$ws.Cells.Item(446, 6).Value2 = 85996
$ws.Cells.Item(446, 7).Value2 = 260
$ws.Cells.Item(447, 6).Value2 = 66696
$ws.Cells.Item(447, 7).Value2 = 270

# --- Append the new day's row (r=448) ---
$ws.Cells.Item(448, 1).Value2 = 44342
$ws.Cells.Item(448, 2).Value2 = 389344
$ws.Cells.Item(448, 3).Value2 = 4977
$ws.Cells.Item(448, 4).Value2 = 168
$ws.Cells.Item(448, 5).Value2 = 12320
$ws.Cells.Item(448, 6).Value2 = 47287
$ws.Cells.Item(448, 7).Value2 = 104
